$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B5 value from "yes" to "no"
$ws.Range("B5").Value = "no"

# Update selection to A5
$ws.Range("A5").Select()
